# Updates Leve profit-calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with refreshed market-board pricing data, as produced by the scheduled
# data-refresh runner. For each affected row, columns H-N
# (currentAveragePrice, currentAveragePriceNQ, currentAveragePriceHQ,
# LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ) are recalculated;
# some rows gain or lose an HQ/NQ profit column depending on whether that
# craft type is applicable.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3336772
$ws.Range("I74").Value = 3336772
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3336772
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3335836
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 3336772
$ws.Range("I77").Value = 3336772
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 16683860
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -16679180
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1908.4762
$ws.Range("I102").Value = 1667.5385
$ws.Range("J102").Value = 2300
$ws.Range("K102").Value = 1667.5385
$ws.Range("L102").Value = 2300
$ws.Range("M102").Value = -45.53850000000011
$ws.Range("N102").Value = -5544

$ws.Range("H110").Value = 705.8889
$ws.Range("I110").Value = 693.2857
$ws.Range("J110").Value = 750
$ws.Range("K110").Value = 693.2857
$ws.Range("L110").Value = 750
$ws.Range("M110").Value = 1351.7143
$ws.Range("N110").Value = -4840

$ws.Range("H122").Value = 789.87177
$ws.Range("I122").Value = 650.15625
$ws.Range("J122").Value = 1428.5714
$ws.Range("K122").Value = 1950.46875
$ws.Range("L122").Value = 4285.7142
$ws.Range("M122").Value = 499.53125
$ws.Range("N122").Value = -9185.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 701840
$ws.Range("J86").Value = 1001771.44
$ws.Range("L86").Value = 1001771.44
$ws.Range("N86").Value = -1004017.44

$ws.Range("H89").Value = 701840
$ws.Range("J89").Value = 1001771.44
$ws.Range("L89").Value = 5008857.199999999
$ws.Range("N89").Value = -5020089.199999999

$ws.Range("H99").Value = 2377.1428
$ws.Range("I99").Value = 2346
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2346
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -848
$ws.Range("N99").Value = -5996

$ws.Range("H107").Value = 1402.4445
$ws.Range("I107").Value = 1202.75
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 1202.75
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 717.25
$ws.Range("N107").Value = -6840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 51.8
$ws.Range("I7").Value = 39.75
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 39.75
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = 73.25
$ws.Range("N7").Value = -326

$ws.Range("H22").Value = 229.66667
$ws.Range("I22").Value = 252.5
$ws.Range("J22").Value = 211.4
$ws.Range("K22").Value = 252.5
$ws.Range("L22").Value = 211.4
$ws.Range("M22").Value = 97.5
$ws.Range("N22").Value = -911.4

$ws.Range("H31").Value = 4490.3584
$ws.Range("I31").Value = 3650
$ws.Range("K31").Value = 3650
$ws.Range("M31").Value = -3355

$ws.Range("H34").Value = 4490.3584
$ws.Range("I34").Value = 3650
$ws.Range("K34").Value = 3650
$ws.Range("M34").Value = -3448

$ws.Range("H62").Value = 62502324
$ws.Range("I62").Value = 71430800
$ws.Range("K62").Value = 71430800
$ws.Range("M62").Value = -71430176

$ws.Range("H65").Value = 62502324
$ws.Range("I65").Value = 71430800
$ws.Range("K65").Value = 357154000
$ws.Range("M65").Value = -357150880

$ws.Range("H97").Value = 22000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 22000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 22000
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -23982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 2703.6667
$ws.Range("I10").Value = 213.5
$ws.Range("J10").Value = 7684
$ws.Range("K10").Value = 640.5
$ws.Range("L10").Value = 23052
$ws.Range("M10").Value = -501.5
$ws.Range("N10").Value = -23330

$ws.Range("H68").Value = 2694.4482
$ws.Range("I68").Value = 4097.148
$ws.Range("J68").Value = 1472.742
$ws.Range("K68").Value = 12291.444
$ws.Range("L68").Value = 4418.226
$ws.Range("M68").Value = -11480.444
$ws.Range("N68").Value = -6040.226

$ws.Range("H71").Value = 2694.4482
$ws.Range("I71").Value = 4097.148
$ws.Range("J71").Value = 1472.742
$ws.Range("K71").Value = 36874.332
$ws.Range("L71").Value = 13254.678
$ws.Range("M71").Value = -32818.332
$ws.Range("N71").Value = -21366.678

$ws.Range("H80").Value = 2955.5557
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2955.5557
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 8866.667099999999
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -10738.6671

$ws.Range("H83").Value = 2955.5557
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2955.5557
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 26600.0013
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -35960.0013

$ws.Range("H107").Value = 469.46295
$ws.Range("I107").Value = 281.93616
$ws.Range("J107").Value = 1728.5714
$ws.Range("K107").Value = 845.8084799999999
$ws.Range("L107").Value = 5185.7142
$ws.Range("M107").Value = 1074.19152
$ws.Range("N107").Value = -9025.7142

$ws.Range("H132").Value = 1172.075
$ws.Range("I132").Value = 909.5833
$ws.Range("J132").Value = 1565.8125
$ws.Range("K132").Value = 8186.2497
$ws.Range("L132").Value = 14092.3125
$ws.Range("M132").Value = -5656.2497
$ws.Range("N132").Value = -19152.3125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14875
$ws.Range("J80").Value = 50500
$ws.Range("L80").Value = 50500
$ws.Range("N80").Value = -52496

$ws.Range("H83").Value = 14875
$ws.Range("J83").Value = 50500
$ws.Range("L83").Value = 252500
$ws.Range("N83").Value = -262484

$ws.Range("H94").Value = 8000
$ws.Range("J94").Value = 8000
$ws.Range("L94").Value = 8000
$ws.Range("N94").Value = -9352

$ws.Range("H97").Value = 1950
$ws.Range("I97").Value = 1920
$ws.Range("J97").Value = 2100
$ws.Range("K97").Value = 1920
$ws.Range("L97").Value = 2100
$ws.Range("M97").Value = -1424
$ws.Range("N97").Value = -3092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 309.60464
$ws.Range("I55").Value = 228.72
$ws.Range("K55").Value = 228.72
$ws.Range("M55").Value = -55.72

$ws.Range("H93").Value = 1534.2
$ws.Range("I93").Value = 1093.625
$ws.Range("J93").Value = 2037.7142
$ws.Range("K93").Value = 1093.625
$ws.Range("L93").Value = 2037.7142
$ws.Range("M93").Value = 154.375
$ws.Range("N93").Value = -4533.7142

$ws.Range("H100").Value = 42151.92
$ws.Range("I100").Value = 101560.4
$ws.Range("J100").Value = 2546.2666
$ws.Range("K100").Value = 101560.4
$ws.Range("L100").Value = 2546.2666
$ws.Range("M100").Value = -101019.4
$ws.Range("N100").Value = -3628.2666
